$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking values (e.g. "1.00")
# are stored as literal strings instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '34.512.44'
$ws.Range("E2").Value = '  +1.66%  '

# Row 3
$ws.Range("D3").Value = '1.838.13'
$ws.Range("E3").Value = '  +3.73%  '

# Row 4
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = '226.34'
$ws.Range("E5").Value = '  +0.63%  '

# Row 6
$ws.Range("D6").Value = '0.554'
$ws.Range("E6").Value = '  +1.44%  '

# Row 7
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.11%  '

# Row 8
$ws.Range("D8").Value = '32.31'
$ws.Range("E8").Value = '  +3.88%  '

# Row 9
$ws.Range("D9").Value = '0.293'
$ws.Range("E9").Value = '  +5.27%  '

# Row 10
$ws.Range("D10").Value = '0.0716'
$ws.Range("E10").Value = '  +9.54%  '

# Row 11
$ws.Range("D11").Value = '0.0934'
$ws.Range("E11").Value = '  +0.60%  '

# Row 12
$ws.Range("D12").Value = '2.109.84'
$ws.Range("E12").Value = '  +4.14%  '

# Row 13
$ws.Range("D13").Value = '1.844.03'
$ws.Range("E13").Value = '  +4.01%  '

# Row 14
$ws.Range("D14").Value = '11.03'
$ws.Range("E14").Value = '  +1.58%  '

# Row 15
$ws.Range("D15").Value = '0.649'
$ws.Range("E15").Value = '  +4.91%  '

# Row 16
$ws.Range("D16").Value = '34.538.01'
$ws.Range("E16").Value = '  +1.74%  '

# Row 17
$ws.Range("D17").Value = '4.35'
$ws.Range("E17").Value = '  +4.06%  '

# Row 18
$ws.Range("D18").Value = '69.77'
$ws.Range("E18").Value = '  +1.99%  '

# Row 19
$ws.Range("D19").Value = '252.75'
$ws.Range("E19").Value = '  +0.82%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0802'
$ws.Range("E20").Value = '  +9.38%  '

# Row 21
$ws.Range("D21").Value = '11.29'
$ws.Range("E21").Value = '  +10.16%  '

# Row 22
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.33%  '

# Row 23
$ws.Range("D23").Value = '4.30'
$ws.Range("E23").Value = '  +3.19%  '

# Row 24
$ws.Range("E24").Value = '  +1.36%  '

# Row 25
$ws.Range("D25").Value = '161.83'
$ws.Range("E25").Value = '  +3.98%  '

# Row 26
$ws.Range("D26").Value = '16.81'
$ws.Range("E26").Value = '  +3.09%  '

# Row 27
$ws.Range("D27").Value = '7.25'
$ws.Range("E27").Value = '  +4.12%  '

# Row 28
$ws.Range("E28").Value = '  +2.10%  '

# Row 29
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.13%  '

# Row 30
$ws.Range("D30").Value = '0.0536'
$ws.Range("E30").Value = '  +5.46%  '

# Row 31
$ws.Range("D31").Value = '3.81'
$ws.Range("E31").Value = '  +1.91%  '

# Row 32
$ws.Range("E32").Value = '  +2.14%  '

# Row 33
$ws.Range("D33").Value = '516.73'
$ws.Range("E33").Value = '  +889.82%  '

# Row 34
$ws.Range("D34").Value = '3.64'
$ws.Range("E34").Value = '  +2.75%  '

# Row 35
$ws.Range("E35").Value = '  +6.58%  '

# Row 36
$ws.Range("D36").Value = '1.457.46'
$ws.Range("E36").Value = '  +0.84%  '

# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.08'
$ws.Range("E37").Value = '  +2.79%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '0.655'
$ws.Range("E38").Value = '  +5.48%  '

# Row 39
$ws.Range("D39").Value = '0.0194'
$ws.Range("E39").Value = '  +4.69%  '

# Row 40
$ws.Range("D40").Value = '0.978'
$ws.Range("E40").Value = '  +11.04%  '

# Row 41
$ws.Range("D41").Value = '82.91'
$ws.Range("E41").Value = '  +0.88%  '

# Row 42
$ws.Range("D42").Value = '2.79'
$ws.Range("E42").Value = '  -1.73%  '

# Row 43
$ws.Range("D43").Value = '2.37'
$ws.Range("E43").Value = '  +1.19%  '

# Row 44
$ws.Range("E44").Value = '  +5.70%  '

# Row 45
$ws.Range("D45").Value = '6.15'
$ws.Range("E45").Value = '  +7.76%  '

# Row 46
$ws.Range("D46").Value = '2.005.48'
$ws.Range("E46").Value = '  +4.24%  '

# Row 47
$ws.Range("E47").Value = '  +1.03%  '

# Row 48
$ws.Range("E48").Value = '  -1.84%  '

# Row 49
$ws.Range("D49").Value = '12.28'
$ws.Range("E49").Value = '  +4.11%  '

# Row 50
$ws.Range("D50").Value = '106.37'
$ws.Range("E50").Value = '  +9.91%  '

# Row 51
$ws.Range("E51").Value = '  +0.16%  '

# Restore the default "Normal" style on column D so no residual text-format
# styling is left applied to the cells (matches original styling).
$ws.Range("D2:D51").Style = "Normal"
